$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing values (B2, C2, B3, C3, D3) ---
$ws.Range("B2").Value = 227.93770000000001
$ws.Range("C2").Value = 109.8947
$ws.Range("B3").Value = 217.3869
$ws.Range("C3").Value = 103.6095
$ws.Range("D3").Value = 9.84

# Apply custom number format to G2:G3 (13 decimal places)
$ws.Range("G2:G3").NumberFormat = "0.0000000000000"

# --- Won / lost samples ---
$ws.Range("B8").Value = 23

# Apply Percent (0.00%) formatting to score/lower bound/upper bound
$ws.Range("B12:B14").Style = "Percent"
$ws.Range("B12:B14").NumberFormat = "0.00%"

# --- Score section ---
$ws.Range("B17").Value = 49.19
$ws.Range("B18").Value = 0.3

# --- New rows 21-23 ---
$ws.Range("A21").Value = "Rank"
$ws.Range("B21").Value = 14

$ws.Range("A22").Value = "Total"
$ws.Range("B22").Value = 2313

$ws.Range("A23").Value = "Percentile"
$ws.Range("B23").Formula = "=B21/B22"
$ws.Range("B23").Style = "Percent"
$ws.Range("B23").NumberFormat = "0.00%"

# --- Column width for column G ---
$ws.Columns.Item(7).ColumnWidth = 21.7109375

# --- Selection / view changes ---
$ws.Range("I16").Select()

$wb.Save()
